# C5-PowerPoint.pptx edit: restore the default "Office Theme" colour palette
# as the presentation's active theme (was the custom "Integral" palette),
# and switch the sources-of-finance table on slide 6 to a different
# built-in table style.

$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Theme colours: Integral -> Office
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in PpThemeColorSchemeIndex
#    order). The deck's single theme part backs both the slide master and
#    the notes master, so this repaints the whole design in one go.
# ------------------------------------------------------------------
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72

# ------------------------------------------------------------------
# 2) Slide 6 table: switch its table style (tableStyleId) away from the
#    workbook default style onto the explicit GUID style.
# ------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{07E19604-548A-4C65-97DF-58FB4B19E270}")
    }
}
